# Trade #102 closed at 2026-02-17 21:27:56 - unknown UNKNOWN +0.000%
# New trade #163 opened at 2026-02-17 21:27:50
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet - roll-up metrics after closing trade #130 (MarketMaking,
# row 131 on "All Trades" / row 98 on "MarketMaking") and opening trade #163.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.36
$wsSummary.Range("B4").Value = 1.15
$wsSummary.Range("B6").Value = 130
$wsSummary.Range("B8").Value = 51
$wsSummary.Range("B9").Value = 42.31

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.36
$wsStatus.Range("D5").Value = 97
$wsStatus.Range("E5").Value = 1.04
$wsStatus.Range("F5").Value = 1.36
$wsStatus.Range("G5").Value = 42.27

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Trade #130 (row 131) transitions from OPEN to CLOSED (early exit)
$wsAll.Range("G131").Value = 0.06
$wsAll.Range("H131").Value = "CLOSED"
$wsAll.Range("I131").Value = -14.2857
$wsAll.Range("J131").Value = -0.01
$wsAll.Range("K131").Value = 101.36
$wsAll.Range("L131").Value = "early_exit"
$wsAll.Range("M131").Value = 0.11

# New trade #163 appended as row 164
$wsAll.Range("A164").Value = 163
$wsAll.Range("B164").Value = "'2026-02-17"
$wsAll.Range("C164").Value = "21:27:50"
$wsAll.Range("D164").Value = "MarketMaking"
$wsAll.Range("E164").Value = "UP"
$wsAll.Range("F164").Value = 0.07000000000000001
$wsAll.Range("G164").Value = "'"
$wsAll.Range("H164").Value = "OPEN"
$wsAll.Range("I164").Value = 0
$wsAll.Range("J164").Value = 0
$wsAll.Range("K164").Value = 101.3741758035408
$wsAll.Range("L164").Value = "'"
$wsAll.Range("M164").Value = 0
$wsAll.Range("N164").Value = 0
$wsAll.Range("O164").Value = 0
$wsAll.Range("P164").Value = 0.6
$wsAll.Range("Q164").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Trade #130 (row 98) transitions from OPEN to CLOSED (early exit)
$wsMM.Range("G98").Value = 0.06
$wsMM.Range("H98").Value = "CLOSED"
$wsMM.Range("I98").Value = -14.2857
$wsMM.Range("J98").Value = -0.01
$wsMM.Range("K98").Value = 101.36
$wsMM.Range("P98").Value = "early_exit"
$wsMM.Range("Q98").Value = 0.11

# New trade #163 appended as row 131
$wsMM.Range("A131").Value = 163
$wsMM.Range("B131").Value = "'2026-02-17"
$wsMM.Range("C131").Value = "21:27:50"
$wsMM.Range("D131").Value = "MarketMaking"
$wsMM.Range("E131").Value = "UP"
$wsMM.Range("F131").Value = 0.07000000000000001
$wsMM.Range("G131").Value = "'"
$wsMM.Range("H131").Value = "OPEN"
$wsMM.Range("I131").Value = 0
$wsMM.Range("J131").Value = 0
$wsMM.Range("K131").Value = 101.3741758035408
$wsMM.Range("L131").Value = 0
$wsMM.Range("M131").Value = 0
$wsMM.Range("N131").Value = 0.6
$wsMM.Range("O131").Value = "Normal spread capture: 19600 bps"
$wsMM.Range("P131").Value = "'"
$wsMM.Range("Q131").Value = 0
